# This script rearranges the data of the three worksheets:
#  - Sheet1 and Sheet2 get the styled "name / gender / age" table that used
#    to live on Sheet3 (keeping Sheet3's cell styles/number formats), but
#    filled in with new people data.
#  - Sheet3 gets the plain "hello / world / good / game" 2x2 table that used
#    to live on Sheet1 / Sheet2.
#  - Selections / the active sheet tab are updated to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Step 1: propagate Sheet3's styled A1:D5 template (with its cell
# styles) onto Sheet1 and Sheet2, before Sheet3 itself gets rebuilt. ---
$ws3.Range("A1:D5").Copy($ws1.Range("A1"))
$ws3.Range("A1:D5").Copy($ws2.Range("A1"))

# --- Step 2: fill Sheet1 with the new table contents. ---
$ws1.Range("A1").Value = "姓名"
$ws1.Range("B1").Value = "性别"
$ws1.Range("C1").Value = "年龄"

$ws1.Range("A2").Value = "张三"
$ws1.Range("B2").Value = "男"
$ws1.Range("C2").Value = 21

$ws1.Range("A3").Value = "李四"
$ws1.Range("B3").Value = "男"
$ws1.Range("C3").Value = 22

$ws1.Range("A4").Value = "李华"
$ws1.Range("B4").Value = "男"
$ws1.Range("C4").Value = 23

$ws1.Range("A5").Value = "王丽"
$ws1.Range("B5").Value = "女"
$ws1.Range("C5").Value = 22

# --- Step 3: fill Sheet2 with the same table contents. ---
$ws2.Range("A1").Value = "姓名"
$ws2.Range("B1").Value = "性别"
$ws2.Range("C1").Value = "年龄"

$ws2.Range("A2").Value = "张三"
$ws2.Range("B2").Value = "男"
$ws2.Range("C2").Value = 21

$ws2.Range("A3").Value = "李四"
$ws2.Range("B3").Value = "男"
$ws2.Range("C3").Value = 22

$ws2.Range("A4").Value = "李华"
$ws2.Range("B4").Value = "男"
$ws2.Range("C4").Value = 23

$ws2.Range("A5").Value = "王丽"
$ws2.Range("B5").Value = "女"
$ws2.Range("C5").Value = 22

# --- Step 4: row heights for the new tables (title row taller). ---
$ws1.Rows.Item(1).RowHeight = 15.6
$ws1.Rows("2:5").RowHeight = 14.4

$ws2.Rows.Item(1).RowHeight = 15.6
$ws2.Rows("2:5").RowHeight = 14.4

# --- Step 5: rebuild Sheet3 as the plain 2x2 "hello/world/good/game" grid
# that used to be on Sheet1 / Sheet2. ---
$ws3.Rows("1:5").Delete()
$ws3.Range("A1").Value = "hello"
$ws3.Range("B1").Value = "world"
$ws3.Range("A2").Value = "good"
$ws3.Range("B2").Value = "game"

# --- Step 6: selections / active sheet. Sheet1 becomes the active tab. ---
$ws2.Select() | Out-Null
$ws2.Range("D13").Select() | Out-Null

$ws3.Select() | Out-Null
$ws3.Range("C2").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("D10").Select() | Out-Null
